$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Treatment query (B5): remove the redundant CONCAT() wrapper ---
$oldTreatmentQuery = $ws.Range("B5").Text
$newTreatmentQuery = $oldTreatmentQuery.Replace(
    "    CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"",",
    "    REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","
)
$ws.Range("B5").Value = $newTreatmentQuery

# Re-apply the cell's wrap/font formatting so it keeps the same look as the
# other query cells in column B.
$ws.Range("B5").WrapText = $true
$ws.Range("B5").Font.Name = "Calibri"
$ws.Range("B5").Font.Size = 12

# --- Leave the cursor/viewport on the fixed cell, as the author did ---
$ws.Range("B5").Select()
$excel.ActiveWindow.ScrollRow = 5
